$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 412
$ws.Cells.Item(2, 9).Value = 468
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 468
$ws.Cells.Item(2, 12).Value = 300
$ws.Cells.Item(2, 13).Value = -355
$ws.Cells.Item(2, 14).Value = -526

# Row 9
$ws.Cells.Item(9, 8).Value = 100000100
$ws.Cells.Item(9, 9).Value = 166666750
$ws.Cells.Item(9, 10).Value = 115.5
$ws.Cells.Item(9, 11).Value = 166666750
$ws.Cells.Item(9, 12).Value = 115.5
$ws.Cells.Item(9, 13).Value = -166666581
$ws.Cells.Item(9, 14).Value = -453.5

# Row 40
$ws.Cells.Item(40, 8).Value = 1926.0344
$ws.Cells.Item(40, 9).Value = 1919.5454
$ws.Cells.Item(40, 10).Value = 1946.4286
$ws.Cells.Item(40, 11).Value = 1919.5454
$ws.Cells.Item(40, 12).Value = 1946.4286
$ws.Cells.Item(40, 13).Value = -1744.5454
$ws.Cells.Item(40, 14).Value = -2296.4286

# Row 51
$ws.Cells.Item(51, 8).Value = 5850
$ws.Cells.Item(51, 9).Value = 1950
$ws.Cells.Item(51, 10).Value = 9750
$ws.Cells.Item(51, 11).Value = 1950
$ws.Cells.Item(51, 12).Value = 9750
$ws.Cells.Item(51, 13).Value = -1466
$ws.Cells.Item(51, 14).Value = -10718

# Row 138
$ws.Cells.Item(138, 8).Value = 3881.0952
$ws.Cells.Item(138, 9).Value = 809.0857
$ws.Cells.Item(138, 10).Value = 7721.107
$ws.Cells.Item(138, 11).Value = 2427.2571
$ws.Cells.Item(138, 12).Value = 23163.321
$ws.Cells.Item(138, 13).Value = 2712.7429
$ws.Cells.Item(138, 14).Value = -33443.321

$ws = $wb.Worksheets.Item("ARM")
# Row 111
$ws.Cells.Item(111, 8).Value = 11644
$ws.Cells.Item(111, 10).Value = 11644
$ws.Cells.Item(111, 12).Value = 11644
$ws.Cells.Item(111, 14).Value = -19824

$ws = $wb.Worksheets.Item("BSM")
# Row 74
$ws.Cells.Item(74, 8).Value = 38261.168
$ws.Cells.Item(74, 10).Value = 45394.2
$ws.Cells.Item(74, 12).Value = 45394.2
$ws.Cells.Item(74, 14).Value = -47266.2

# Row 77
$ws.Cells.Item(77, 8).Value = 38261.168
$ws.Cells.Item(77, 10).Value = 45394.2
$ws.Cells.Item(77, 12).Value = 136182.6
$ws.Cells.Item(77, 14).Value = -145542.6

# Row 107
$ws.Cells.Item(107, 8).Value = 1339.9
$ws.Cells.Item(107, 9).Value = 1339.9
$ws.Cells.Item(107, 11).Value = 1339.9
$ws.Cells.Item(107, 13).Value = 580.0999999999999

# Row 134
$ws.Cells.Item(134, 8).Value = 4300.1284
$ws.Cells.Item(134, 9).Value = 4667.9033
$ws.Cells.Item(134, 10).Value = 2875
$ws.Cells.Item(134, 11).Value = 14003.7099
$ws.Cells.Item(134, 12).Value = 8625
$ws.Cells.Item(134, 13).Value = -11468.7099
$ws.Cells.Item(134, 14).Value = -13695

# Row 139
$ws.Cells.Item(139, 8).Value = 50567.25
$ws.Cells.Item(139, 9).Value = 20709
$ws.Cells.Item(139, 10).Value = 60520
$ws.Cells.Item(139, 11).Value = 20709
$ws.Cells.Item(139, 12).Value = 60520
$ws.Cells.Item(139, 13).Value = -15569
$ws.Cells.Item(139, 14).Value = -70800

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 3216.5476
$ws.Cells.Item(31, 9).Value = 1468.7142
$ws.Cells.Item(31, 10).Value = 4964.381
$ws.Cells.Item(31, 11).Value = 1468.7142
$ws.Cells.Item(31, 12).Value = 4964.381
$ws.Cells.Item(31, 13).Value = -1173.7142
$ws.Cells.Item(31, 14).Value = -5554.381

# Row 34
$ws.Cells.Item(34, 8).Value = 3216.5476
$ws.Cells.Item(34, 9).Value = 1468.7142
$ws.Cells.Item(34, 10).Value = 4964.381
$ws.Cells.Item(34, 11).Value = 1468.7142
$ws.Cells.Item(34, 12).Value = 4964.381
$ws.Cells.Item(34, 13).Value = -1266.7142
$ws.Cells.Item(34, 14).Value = -5368.381

# Row 62
$ws.Cells.Item(62, 8).Value = 7233.3335
$ws.Cells.Item(62, 9).Value = 8350
$ws.Cells.Item(62, 11).Value = 8350
$ws.Cells.Item(62, 13).Value = -7726

# Row 65
$ws.Cells.Item(65, 8).Value = 7233.3335
$ws.Cells.Item(65, 9).Value = 8350
$ws.Cells.Item(65, 11).Value = 41750
$ws.Cells.Item(65, 13).Value = -38630

# Row 98
$ws.Cells.Item(98, 8).Value = 27593.334
$ws.Cells.Item(98, 10).Value = 27593.334
$ws.Cells.Item(98, 12).Value = 27593.334
$ws.Cells.Item(98, 14).Value = -32085.334

# Row 132
$ws.Cells.Item(132, 8).Value = 1991.6305
$ws.Cells.Item(132, 9).Value = 1813.7941
$ws.Cells.Item(132, 10).Value = 2495.5
$ws.Cells.Item(132, 11).Value = 5441.3823
$ws.Cells.Item(132, 12).Value = 7486.5
$ws.Cells.Item(132, 13).Value = -2911.3823
$ws.Cells.Item(132, 14).Value = -12546.5

# Row 134
$ws.Cells.Item(134, 8).Value = 2527.8333
$ws.Cells.Item(134, 9).Value = 2709.4
$ws.Cells.Item(134, 10).Value = 1620
$ws.Cells.Item(134, 11).Value = 8128.200000000001
$ws.Cells.Item(134, 12).Value = 4860
$ws.Cells.Item(134, 13).Value = -5593.200000000001
$ws.Cells.Item(134, 14).Value = -9930

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Cells.Item(32, 8).Value = 1311.25
$ws.Cells.Item(32, 10).Value = 1355.7142
$ws.Cells.Item(32, 12).Value = 4067.1426
$ws.Cells.Item(32, 14).Value = -4633.142599999999

# Row 114
$ws.Cells.Item(114, 8).Value = 7108.2354
$ws.Cells.Item(114, 9).Value = 518.5
$ws.Cells.Item(114, 10).Value = 12965.777
$ws.Cells.Item(114, 11).Value = 1555.5
$ws.Cells.Item(114, 12).Value = 38897.331
$ws.Cells.Item(114, 13).Value = 1698.5
$ws.Cells.Item(114, 14).Value = -45405.331

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 30425672
$ws.Cells.Item(122, 9).Value = 39439628
$ws.Cells.Item(122, 11).Value = 118318884
$ws.Cells.Item(122, 13).Value = -118316434

# Row 132
$ws.Cells.Item(132, 8).Value = 2737.0981
$ws.Cells.Item(132, 9).Value = 2722.0605
$ws.Cells.Item(132, 10).Value = 2764.6667
$ws.Cells.Item(132, 11).Value = 8166.181500000001
$ws.Cells.Item(132, 12).Value = 8294.000100000001
$ws.Cells.Item(132, 13).Value = -5636.181500000001
$ws.Cells.Item(132, 14).Value = -13354.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 8).Value = 2708.2727
$ws.Cells.Item(61, 9).Value = 2040.4667
$ws.Cells.Item(61, 10).Value = 4139.2856
$ws.Cells.Item(61, 11).Value = 2040.4667
$ws.Cells.Item(61, 12).Value = 4139.2856
$ws.Cells.Item(61, 13).Value = -1838.4667
$ws.Cells.Item(61, 14).Value = -4543.2856

# Row 113
$ws.Cells.Item(113, 8).Value = 2708.2727
$ws.Cells.Item(113, 9).Value = 2040.4667
$ws.Cells.Item(113, 10).Value = 4139.2856
$ws.Cells.Item(113, 11).Value = 2040.4667
$ws.Cells.Item(113, 12).Value = 4139.2856
$ws.Cells.Item(113, 13).Value = 129.5333000000001
$ws.Cells.Item(113, 14).Value = -8479.285599999999

# Row 122
$ws.Cells.Item(122, 8).Value = 10182984
$ws.Cells.Item(122, 9).Value = 17863718
$ws.Cells.Item(122, 10).Value = 2502250
$ws.Cells.Item(122, 11).Value = 53591154
$ws.Cells.Item(122, 12).Value = 7506750
$ws.Cells.Item(122, 13).Value = -53588704
$ws.Cells.Item(122, 14).Value = -7511650

# Row 127
$ws.Cells.Item(127, 8).Value = 40399.6
$ws.Cells.Item(127, 10).Value = 40399.6
$ws.Cells.Item(127, 12).Value = 40399.6
$ws.Cells.Item(127, 14).Value = -50319.6

# Row 132
$ws.Cells.Item(132, 8).Value = 7771981.5
$ws.Cells.Item(132, 9).Value = 11462622
$ws.Cells.Item(132, 10).Value = 2211.3157
$ws.Cells.Item(132, 11).Value = 34387866
$ws.Cells.Item(132, 12).Value = 6633.9471
$ws.Cells.Item(132, 13).Value = -34385336
$ws.Cells.Item(132, 14).Value = -11693.9471

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).ClearContents()

# Row 81
$ws.Cells.Item(81, 8).Value = 2493.4
$ws.Cells.Item(81, 9).Value = 2277.889
$ws.Cells.Item(81, 10).Value = 2816.6667
$ws.Cells.Item(81, 11).Value = 4555.778
$ws.Cells.Item(81, 12).Value = 5633.3334
$ws.Cells.Item(81, 13).Value = -3494.778
$ws.Cells.Item(81, 14).Value = -7755.3334

# Row 84
$ws.Cells.Item(84, 8).Value = 2493.4
$ws.Cells.Item(84, 9).Value = 2277.889
$ws.Cells.Item(84, 10).Value = 2816.6667
$ws.Cells.Item(84, 11).Value = 22778.89
$ws.Cells.Item(84, 12).Value = 28166.667
$ws.Cells.Item(84, 13).Value = -17474.89
$ws.Cells.Item(84, 14).Value = -38774.667

# Row 122
$ws.Cells.Item(122, 8).Value = 1457
$ws.Cells.Item(122, 9).Value = 1479.8
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 11).Value = 4439.4
$ws.Cells.Item(122, 12).Value = 4200
$ws.Cells.Item(122, 13).Value = -1989.4
$ws.Cells.Item(122, 14).Value = -9100

# Row 132
$ws.Cells.Item(132, 8).Value = 1664.0605
$ws.Cells.Item(132, 9).Value = 1074.579
$ws.Cells.Item(132, 11).Value = 3223.737
$ws.Cells.Item(132, 13).Value = -693.7370000000001
